$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1: 94.53 -> 0M
$t.Cell(1,1).Range.Text = "0M"

# Row 2: 69.27 -> 0M
$t.Cell(2,1).Range.Text = "0M"

# Row 3: 1267 -> 0M
$t.Cell(3,1).Range.Text = "0M"

# Row 4: 3317 -> 3562
$t.Cell(4,1).Range.Text = "3562"

# Row 7: 0.04998 -> 0.04835
$t.Cell(7,1).Range.Text = "0.04835"

# Row 8: 0.03561 -> 0.03480
$t.Cell(8,1).Range.Text = "0.03480"

# Row 12: 61.92564 -> 69.26598
$t.Cell(12,1).Range.Text = "69.26598"

# Row 44: collapse multi-run tab-separated content down to "94.53"
$t.Cell(44,1).Range.Text = "94.53"

# Row 45: collapse multi-run tab-separated content down to "69.27"
$t.Cell(45,1).Range.Text = "69.27"

# Row 46: collapse multi-run tab-separated content down to "1267"
$t.Cell(46,1).Range.Text = "1267"
